# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets. Both sheets carry identical data; F2/F6/F7/F9 each bump by a
# small amount.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 126
    $ws.Range("F6").Value = 435
    $ws.Range("F7").Value = 149
    $ws.Range("F9").Value = 553
}
